$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column ("Price") values are plain decimal-looking strings (e.g. "588.07").
# Assigning such a string straight to Range.Value on a General-formatted cell
# gets auto-coerced to a number by Excel, which would change both the stored
# type and the literal text (e.g. trailing zeros). Force text entry by
# temporarily marking the cell as Text, assigning, then clearing the format
# back to General so the cell style is left exactly as it was.
function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = '@'
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue 'D2' '63.556.43'
$ws.Range('E2').Value = '  -1.51%  '
Set-TextValue 'D3' '3.065.78'
$ws.Range('E3').Value = '  -4.09%  '
$ws.Range('E4').Value = '  -0.22%  '
Set-TextValue 'D5' '588.07'
$ws.Range('E5').Value = '  -1.11%  '
Set-TextValue 'D6' '154.02'
$ws.Range('E6').Value = '  +3.44%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -0.10%  '
Set-TextValue 'D9' '3.064.74'
$ws.Range('E9').Value = '  -3.86%  '
$ws.Range('E10').Value = '  -4.83%  '
Set-TextValue 'D11' '5.82'
$ws.Range('E11').Value = '  -2.80%  '
$ws.Range('E12').Value = '  -3.38%  '
Set-TextValue 'D13' '36.72'
$ws.Range('E13').Value = '  -2.98%  '
$ws.Range('E14').Value = '  -5.06%  '
$ws.Range('E15').Value = '  -2.43%  '
Set-TextValue 'D16' '3.572.00'
$ws.Range('E16').Value = '  -4.10%  '
Set-TextValue 'D17' '63.547.78'
$ws.Range('E17').Value = '  -1.15%  '
$ws.Range('E18').Value = '  -3.77%  '
Set-TextValue 'D19' '3.064.47'
$ws.Range('E19').Value = '  -3.93%  '
Set-TextValue 'D20' '468.90'
$ws.Range('E20').Value = '  -1.41%  '
Set-TextValue 'D21' '14.24'
$ws.Range('E21').Value = '  -2.95%  '
Set-TextValue 'D22' '0.700'
$ws.Range('E22').Value = '  -5.79%  '
$ws.Range('E23').Value = '  -3.68%  '
Set-TextValue 'D24' '2.42'
$ws.Range('E24').Value = '  -0.88%  '
Set-TextValue 'D25' '80.29'
$ws.Range('E25').Value = '  -2.03%  '
Set-TextValue 'D26' '12.70'
$ws.Range('E26').Value = '  -4.34%  '
Set-TextValue 'D27' '10.40'
$ws.Range('E27').Value = '  +2.76%  '
$ws.Range('E28').Value = '  -0.33%  '
Set-TextValue 'D29' '7.33'
$ws.Range('E29').Value = '  +0.51%  '
$ws.Range('E30').Value = '  -0.28%  '
$ws.Range('E31').Value = '  -3.57%  '
Set-TextValue 'D32' '2.14'
$ws.Range('E32').Value = '  -6.17%  '
Set-TextValue 'D33' '27.03'
$ws.Range('E33').Value = '  -5.43%  '
$ws.Range('E34').Value = '  -6.33%  '
Set-TextValue 'D35' '0.0₃0817'
$ws.Range('E35').Value = '  -5.60%  '
$ws.Range('E36').Value = '  -3.28%  '
Set-TextValue 'D37' '5.96'
$ws.Range('E37').Value = '  -4.66%  '
Set-TextValue 'D38' '3.25'
$ws.Range('E38').Value = '  -3.07%  '
Set-TextValue 'D39' '2.21'
$ws.Range('E39').Value = '  -5.27%  '
Set-TextValue 'D40' '50.40'
$ws.Range('E40').Value = '  -2.26%  '
Set-TextValue 'D41' '9.14'
$ws.Range('E41').Value = '  -3.37%  '
Set-TextValue 'D42' '436.59'
$ws.Range('E42').Value = '  -7.21%  '
Set-TextValue 'D43' '0.286'
$ws.Range('E43').Value = '  -3.89%  '
Set-TextValue 'D44' '40.48'
$ws.Range('E44').Value = '  +2.03%  '
Set-TextValue 'D45' '0.111'
$ws.Range('E45').Value = '  +1.20%  '
$ws.Range('E46').Value = '  -5.69%  '
Set-TextValue 'D47' '2.788.96'
$ws.Range('E47').Value = '  -5.02%  '
Set-TextValue 'D48' '129.02'
$ws.Range('E48').Value = '  -3.20%  '
Set-TextValue 'D50' '24.86'
$ws.Range('E50').Value = '  +0.74%  '
Set-TextValue 'D51' '2.21'
$ws.Range('E51').Value = '  -3.41%  '